{"js": "// \"Trying to get Wallrun on Curved Surface\"\n//\n// In the \"Wall Run\" row of the status table, the notes cell ends with the\n// paragraph \"Maintaining speed on a curved surface.\" Remove the four\n// paragraphs that used to follow it (a blank paragraph, the \"Player\n// movement when face is against the wall...\" note, another blank\n// paragraph, and the \"If Player is facing the wall...\" note), leaving the\n// cell ending at \"Maintaining speed on a curved surface.\"\n\nconst body = context.document.body;\n\n// Locate the anchor paragraph via search so the edit is resilient to the\n// exact paragraph index in the document.\nconst searchResults = body.search(\"Maintaining speed on a curved surface.\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find anchor paragraph \"Maintaining speed on a curved surface.\"');\n}\n\nconst anchorParagraph = searchResults.items[0].paragraphs.getFirst();\nawait context.sync();\n\n// Walk forward from the anchor and collect the next four paragraphs, which\n// are the ones that need to be removed.\nconst toDelete = [];\nlet current = anchorParagraph;\nfor (let i = 0; i < 4; i++) {\n  current = current.getNext();\n  toDelete.push(current);\n}\ntoDelete.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nconst expectedTexts = [\n  \"\",\n  \"Player movement when face is against the wall. (Not the right or left side of the player) \",\n  \"\",\n  \"If Player is facing the wall have left and right movement move player properly along the wall.\",\n];\n\nfor (let i = 0; i < toDelete.length; i++) {\n  const actual = toDelete[i].text.trim();\n  const expected = expectedTexts[i].trim();\n  if (actual !== expected) {\n    throw new Error(\n      `Unexpected paragraph content while removing curved-surface wall-run notes. Expected \"${expected}\" but found \"${actual}\".`\n    );\n  }\n}\n\n// Delete from the last paragraph back to the first so earlier deletions\n// don't shift the identity of paragraphs still queued for removal.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\nawait context.sync();\n", "ps1": "# \"Trying to get Wallrun on Curved Surface\"\n#\n# In the \"Wall Run\" row of the status table, the notes cell ends with the\n# paragraph \"Maintaining speed on a curved surface.\" Remove the four\n# paragraphs that used to follow it (a blank paragraph, the \"Player\n# movement when face is against the wall...\" note, another blank\n# paragraph, and the \"If Player is facing the wall...\" note), leaving the\n# cell ending at \"Maintaining speed on a curved surface.\"\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph via Find so the edit is resilient to the\n# exact paragraph index in the document.\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"Maintaining speed on a curved surface.\")\nif (-not $found) {\n  throw 'Could not find anchor paragraph \"Maintaining speed on a curved surface.\"'\n}\n\n$paras = $d.Paragraphs\n$anchorIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n  if ($paras.Item($i).Range.Start -eq $searchRange.Start) {\n    $anchorIndex = $i\n    break\n  }\n}\nif ($anchorIndex -eq -1) {\n  throw \"Could not resolve anchor paragraph index\"\n}\n\n# The four paragraphs immediately following the anchor are the ones that\n# need to be removed. (The last one also carries the table cell's\n# end-of-cell mark, so compare after trimming paragraph/cell-mark chars.)\n$expectedTexts = @(\n  \"\",\n  \"Player movement when face is against the wall. (Not the right or left side of the player) \",\n  \"\",\n  \"If Player is facing the wall have left and right movement move player properly along the wall.\"\n)\n\n$indicesToDelete = @($anchorIndex + 1, $anchorIndex + 2, $anchorIndex + 3, $anchorIndex + 4)\n\nfor ($j = 0; $j -lt $indicesToDelete.Length; $j++) {\n  $idx = $indicesToDelete[$j]\n  $actual = $d.Paragraphs.Item($idx).Range.Text.TrimEnd(\"`r\", \"`a\")\n  $expected = $expectedTexts[$j]\n  if ($actual -ne $expected) {\n    throw \"Unexpected paragraph content while removing curved-surface wall-run notes at index $idx. Expected [$expected] but found [$actual].\"\n  }\n}\n\n# Delete from the last paragraph back to the first so earlier deletions\n# don't shift the identity of paragraphs still queued for removal.\nfor ($j = $indicesToDelete.Length - 1; $j -ge 0; $j--) {\n  $idx = $indicesToDelete[$j]\n  $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
